$wb = $excel.ActiveWorkbook

# New header labels shared by both the "UseCase" and "UseCase1" sheets.
$headers = @{
    "D1" = "use_case_category"
    "E1" = "known_limitations"
    "F1" = "relevance_to_dgps"
    "G1" = "data_types"
    "H1" = "data_substrates"
    "I1" = "standards_and_tools_for_dgp_use"
    "J1" = "alternative_standards_and_tools"
    "K1" = "enables"
    "L1" = "involved_in_experimental_design"
    "M1" = "involved_in_metadata_management"
    "N1" = "involved_in_quality_control"
    "O1" = "xrefs"
}

foreach ($sheetName in @("UseCase", "UseCase1")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($addr in $headers.Keys) {
        $ws.Range($addr).Value = $headers[$addr]
    }
}

# The "UseCase" sheet additionally carries list-based data validation rules.
# Replace the old "vital_status" (column G) validation with new rules that
# target the new "use_case_category" (D) and "relevance_to_dgps" (F) columns.
$wsUseCase = $wb.Worksheets.Item("UseCase")

$wsUseCase.Range("G2:G1048576").Validation.Delete()

$wsUseCase.Range("D2:D1048576").Validation.Add(3, 1, 1, """Acquisition,Integration,Standardization,Modeling,Application,Assessment""")
$wsUseCase.Range("F2:F1048576").Validation.Add(3, 1, 1, """AI-READI,CHoRUS,CM4AI,Voice""")
